$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1200032159686089, 0.9604939222335815, 0.008076382800936699, 0.9986888766288757),
    @(0.03782502189278603, 0.987867534160614, 0.006330511532723904, 0.9986888766288757),
    @(0.02291683107614517, 0.9947972893714905, 0.003965023905038834, 0.9988905787467957),
    @(0.01646891422569752, 0.9954233169555664, 0.002476252149790525, 0.9989914298057556),
    @(0.0110191497951746, 0.9973015189170837, 0.002119767246767879, 0.9990922808647156),
    @(0.01016345806419849, 0.9972151517868042, 0.002026139525696635, 0.9990922808647156),
    @(0.008427067659795284, 0.997603714466095, 0.00176313822157681, 0.9990922808647156),
    @(0.008846169337630272, 0.9974526166915894, 0.001878790906630456, 0.9990922808647156),
    @(0.008235546760261059, 0.9976900815963745, 0.001599812181666493, 0.9990922808647156),
    @(0.007803684566169977, 0.9977116584777832, 0.001498537021689117, 0.9990922808647156),
    @(0.00871565192937851, 0.9972583055496216, 0.001423503272235394, 0.9990922808647156),
    @(0.008151769638061523, 0.9975174069404602, 0.001271615619771183, 0.9990922808647156),
    @(0.008394578471779823, 0.9972367286682129, 0.001427231589332223, 0.9990922808647156),
    @(0.008263843134045601, 0.9972367286682129, 0.001215763739310205, 0.9990922808647156),
    @(0.007269247900694609, 0.9976469278335571, 0.001106060342863202, 0.9990922808647156),
    @(0.008388400077819824, 0.9974957704544067, 0.001095374813303351, 0.9990922808647156),
    @(0.006611845456063747, 0.9977764487266541, 0.001058319234289229, 0.9990922808647156),
    @(0.007537002675235271, 0.9973878264427185, 0.00121416908223182, 0.9990922808647156),
    @(0.00700663635507226, 0.9975821375846863, 0.001027434598654509, 0.9990922808647156),
    @(0.007252276875078678, 0.9974957704544067, 0.001082345959730446, 0.9990922808647156),
    @(0.006978060584515333, 0.9975605607032776, 0.001330017577856779, 0.9990922808647156),
    @(0.006980582606047392, 0.9975605607032776, 0.0009659050847403705, 0.9990922808647156),
    @(0.007618582807481289, 0.9973015189170837, 0.0009657004266045988, 0.9990922808647156),
    @(0.007244183216243982, 0.9974526166915894, 0.001034227665513754, 0.9990922808647156),
    @(0.007028030697256327, 0.9974526166915894, 0.0009060431038960814, 0.9990922808647156),
    @(0.007348047569394112, 0.9974526166915894, 0.0009090504609048367, 0.9990922808647156),
    @(0.00712881563231349, 0.9974957704544067, 0.001214290852658451, 0.9990922808647156),
    @(0.006449383683502674, 0.9977116584777832, 0.0009924278128892183, 0.9990922808647156),
    @(0.006217965390533209, 0.9977764487266541, 0.0008789451094344258, 0.9990922808647156),
    @(0.006972514092922211, 0.9975389838218689, 0.00113665871322155, 0.9990922808647156),
    @(0.006890335585922003, 0.9975605607032776, 0.0009421516442671418, 0.9990922808647156),
    @(0.007153441663831472, 0.9975174069404602, 0.0008716598385944963, 0.9990922808647156),
    @(0.006508206017315388, 0.9976469278335571, 0.0008188926731236279, 0.9990922808647156),
    @(0.008802797645330429, 0.9972583055496216, 0.001260441727936268, 0.9990922808647156),
    @(0.007221163250505924, 0.9975174069404602, 0.0008302931673824787, 0.9990922808647156),
    @(0.007102751173079014, 0.9975389838218689, 0.0008383350213989615, 0.9990922808647156),
    @(0.006864743307232857, 0.9974957704544067, 0.0008201253949664533, 0.9990922808647156),
    @(0.006658699363470078, 0.9975389838218689, 0.0008612087112851441, 0.9990922808647156),
    @(0.006486013531684875, 0.9976469278335571, 0.000849985342938453, 0.9990922808647156),
    @(0.005958870984613895, 0.9978843927383423, 0.0008277925662696362, 0.9990922808647156),
    @(0.005942939780652523, 0.9977980256080627, 0.0007642018026672304, 0.9990922808647156),
    @(0.007779129780828953, 0.9975821375846863, 0.0007554296171292663, 0.9990922808647156),
    @(0.006939390208572149, 0.9975174069404602, 0.000766658631619066, 0.9990922808647156),
    @(0.006943475920706987, 0.9975389838218689, 0.0007845893851481378, 0.9990922808647156),
    @(0.00606096163392067, 0.9977332353591919, 0.0007758575957268476, 0.9990922808647156),
    @(0.00770706869661808, 0.9972583055496216, 0.0008030487224459648, 0.9990922808647156),
    @(0.007006468717008829, 0.997474193572998, 0.0007973500760272145, 0.9990922808647156),
    @(0.006244299001991749, 0.9976685047149658, 0.0007813083357177675, 0.9990922808647156),
    @(0.007039610296487808, 0.9974310398101807, 0.0008661326137371361, 0.9990922808647156),
    @(0.00634013582020998, 0.9977332353591919, 0.001469014910981059, 0.9990922808647156)
)
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
